# Update odds values in Sheet1 for row 2 and row 5, reflecting refreshed
# FlashScore odds data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 (U. Magdalena vs Millonarios) ---
$ws.Range("G2").Value  = 3.1
$ws.Range("H2").Value  = 3.1
$ws.Range("I2").Value  = 2.38
$ws.Range("K2").Value  = 1.95
$ws.Range("L2").Value  = 3.25
$ws.Range("M2").Value  = 1.1
$ws.Range("N2").Value  = 7
$ws.Range("O2").Value  = 1.44
$ws.Range("P2").Value  = 2.63
$ws.Range("Q2").Value  = 1.85
$ws.Range("R2").Value  = 2
$ws.Range("S2").Value  = 2.4
$ws.Range("T2").Value  = 1.53
$ws.Range("U2").Value  = 4.2
$ws.Range("V2").Value  = 1.22
$ws.Range("W2").Value  = 5
$ws.Range("X2").Value  = 1.17
$ws.Range("Y2").Value  = 1.57
$ws.Range("Z2").Value  = 2.25
$ws.Range("AA2").Value = 2.1
$ws.Range("AB2").Value = 1.67
$ws.Range("AC2").Value = 7.5
$ws.Range("AI2").Value = 7
$ws.Range("AL2").Value = 67
$ws.Range("AM2").Value = 6.5
$ws.Range("AN2").Value = 10
$ws.Range("AO2").Value = 10
$ws.Range("AQ2").Value = 23

# --- Row 5 (Correcaminos vs Dorados de Sinaloa) ---
$ws.Range("G5").Value  = 3.45
$ws.Range("H5").Value  = 3.2
$ws.Range("I5").Value  = 2.05
$ws.Range("J5").Value  = 3.85
$ws.Range("L5").Value  = 2.62
$ws.Range("O5").Value  = 1.23
$ws.Range("P5").Value  = 3.35
$ws.Range("S5").Value  = 1.7
$ws.Range("T5").Value  = 1.91
$ws.Range("W5").Value  = 2.62
$ws.Range("X5").Value  = 1.37
$ws.Range("AB5").Value = 2.12
$ws.Range("AC5").Value = 12
$ws.Range("AD5").Value = 21
$ws.Range("AE5").Value = 11.5
$ws.Range("AF5").Value = 50
$ws.Range("AG5").Value = 29
$ws.Range("AH5").Value = 30
$ws.Range("AK5").Value = 11.75
$ws.Range("AL5").Value = 45
$ws.Range("AM5").Value = 8.5
$ws.Range("AN5").Value = 10.75
$ws.Range("AO5").Value = 8.25
$ws.Range("AP5").Value = 20
$ws.Range("AQ5").Value = 15.5
$ws.Range("AR5").Value = 22

$wb.Save()
